# fix[DE]: Adjust tests for CandidateInstance
#
# The workbook has two worksheets used as test fixtures for the Data Set
# Explorer ("Long Method Smell" and "Large Class Smell"). This edit:
#   1. Renames both sheets to drop the trailing "Smell" suffix.
#   2. Switches the active/selected tab from "Large Class" back to
#      "Long Method" and resets its view/selection to the top of the sheet.

$wb = $excel.ActiveWorkbook

$wsLongMethod = $wb.Worksheets.Item(1)
$wsLargeClass = $wb.Worksheets.Item(2)

$wsLongMethod.Name = "Long Method"
$wsLargeClass.Name = "Large Class"

# Make "Long Method" the active/visible tab again and move the selection
# to E12 (the frozen header pane keeps rows 1-3 pinned, so this scrolls
# the sheet back up near the top).
$wsLongMethod.Activate()
$wsLongMethod.Range("E12").Select()
